$d = $word.ActiveDocument

$d.Content.Find.Execute("12×60=720", $true, $false, $false, $false, $false, $true, 1, $false, "45×55=2475", 2) | Out-Null
$d.Content.Find.Execute("60×20=1200", $true, $false, $false, $false, $false, $true, 1, $false, "41×14=574", 2) | Out-Null
$d.Content.Find.Execute("22×81=1782", $true, $false, $false, $false, $false, $true, 1, $false, "40×81=3240", 2) | Out-Null
$d.Content.Find.Execute("54×87=4698", $true, $false, $false, $false, $false, $true, 1, $false, "89×29=2581", 2) | Out-Null
$d.Content.Find.Execute("60×94=5640", $true, $false, $false, $false, $false, $true, 1, $false, "50×11=550", 2) | Out-Null
$d.Content.Find.Execute("61×15=915", $true, $false, $false, $false, $false, $true, 1, $false, "63×46=2898", 2) | Out-Null
$d.Content.Find.Execute("68×16=1088", $true, $false, $false, $false, $false, $true, 1, $false, "37×52=1924", 2) | Out-Null
$d.Content.Find.Execute("29×11=319", $true, $false, $false, $false, $false, $true, 1, $false, "99×29=2871", 2) | Out-Null
$d.Content.Find.Execute("85×88=7480", $true, $false, $false, $false, $false, $true, 1, $false, "99×75=7425", 2) | Out-Null
$d.Content.Find.Execute("82×89=7298", $true, $false, $false, $false, $false, $true, 1, $false, "16×65=1040", 2) | Out-Null
$d.Content.Find.Execute("88×63=5544", $true, $false, $false, $false, $false, $true, 1, $false, "52×99=5148", 2) | Out-Null
$d.Content.Find.Execute("49×77=3773", $true, $false, $false, $false, $false, $true, 1, $false, "11×97=1067", 2) | Out-Null
$d.Content.Find.Execute("21×66=1386", $true, $false, $false, $false, $false, $true, 1, $false, "97×27=2619", 2) | Out-Null
$d.Content.Find.Execute("15×67=1005", $true, $false, $false, $false, $false, $true, 1, $false, "37×93=3441", 2) | Out-Null
$d.Content.Find.Execute("62×36=2232", $true, $false, $false, $false, $false, $true, 1, $false, "64×50=3200", 2) | Out-Null
$d.Content.Find.Execute("48×99=4752", $true, $false, $false, $false, $false, $true, 1, $false, "28×59=1652", 2) | Out-Null
$d.Content.Find.Execute("12×75=900", $true, $false, $false, $false, $false, $true, 1, $false, "45×46=2070", 2) | Out-Null
$d.Content.Find.Execute("19×79=1501", $true, $false, $false, $false, $false, $true, 1, $false, "35×74=2590", 2) | Out-Null
$d.Content.Find.Execute("86×32=2752", $true, $false, $false, $false, $false, $true, 1, $false, "21×52=1092", 2) | Out-Null
$d.Content.Find.Execute("61×73=4453", $true, $false, $false, $false, $false, $true, 1, $false, "35×96=3360", 2) | Out-Null
$d.Content.Find.Execute("89×85=7565", $true, $false, $false, $false, $false, $true, 1, $false, "57×38=2166", 2) | Out-Null
$d.Content.Find.Execute("82×53=4346", $true, $false, $false, $false, $false, $true, 1, $false, "52×36=1872", 2) | Out-Null
$d.Content.Find.Execute("73×85=6205", $true, $false, $false, $false, $false, $true, 1, $false, "72×60=4320", 2) | Out-Null
$d.Content.Find.Execute("27×22=594", $true, $false, $false, $false, $false, $true, 1, $false, "65×18=1170", 2) | Out-Null
$d.Content.Find.Execute("29×61=1769", $true, $false, $false, $false, $false, $true, 1, $false, "76×55=4180", 2) | Out-Null
